$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the header row (row 1), pushing the
# existing dividend history rows down one, and populate it with the
# newest dividend entry (06/01/2025, 0.010).
$ws.Range("A2:C2").EntireRow.Insert()

# Leading apostrophe forces these to be stored as text (matching the
# existing "03/10/2024" / "0.010" cells, which are plain text strings,
# not real dates/numbers).
$ws.Range("A2").Value = "'06/01/2025"
$ws.Range("B2").Value = "'06/01/2025"
$ws.Range("C2").Value = "'0.010"

# Drop the implicit "quote prefix" formatting Excel applies for the
# leading apostrophe so the new cells keep the sheet's plain default
# style, same as every other cell in this sheet.
$ws.Range("A2:C2").ClearFormats()
